# Apply updated metrics data to Sheet1 (rows 2-26, columns A-I)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("model_12_1_17", 0.7702880159603739, 0.1987834259997685, 0.3737867480225621, 0.4677666912442467, 0.254223495721817, 0.1547020375728607, 0.5813519358634949, 0.3554784059524536),
    @("model_12_1_12", 0.7702981945054015, 0.1996665894790683, 0.3737078230375384, 0.4678501598792255, 0.2542122006416321, 0.1545315086841583, 0.5814251899719238, 0.355422705411911),
    @("model_12_1_19", 0.7703152824586406, 0.1989619078212664, 0.3739664464717619, 0.4679115358198676, 0.2541933059692383, 0.1546675562858582, 0.5811851024627686, 0.3553816676139832),
    @("model_12_1_20", 0.7703187006311639, 0.1982150082930486, 0.3741693050712878, 0.4679298873447452, 0.254189521074295, 0.1548117846250534, 0.5809967517852783, 0.3553694188594818),
    @("model_12_1_23", 0.770333287917822, 0.1993957070667797, 0.3740086516862537, 0.4680055387702095, 0.2541733682155609, 0.1545837968587875, 0.5811458826065063, 0.3553188741207123),
    @("model_12_1_15", 0.7703431768201322, 0.2001708005561335, 0.3739285619974192, 0.4680717698651028, 0.2541624009609222, 0.1544341444969177, 0.581220269203186, 0.3552746772766113),
    @("model_12_1_21", 0.7703525420718753, 0.1987242055308958, 0.3743293730948457, 0.4681125161889585, 0.2541520595550537, 0.1547134518623352, 0.5808481574058533, 0.355247437953949),
    @("model_12_1_24", 0.7703547736059736, 0.1993563993320111, 0.374195741640775, 0.4681218895807548, 0.2541496157646179, 0.1545913964509964, 0.5809721946716309, 0.3552411496639252),
    @("model_12_1_11", 0.7703641322270097, 0.1991967505715344, 0.374369412421717, 0.4682110444385842, 0.254139244556427, 0.1546222269535065, 0.5808110237121582, 0.3551816344261169),
    @("model_12_1_22", 0.7703641853590977, 0.1989926258025301, 0.374361939156395, 0.4681748873317186, 0.2541391849517822, 0.1546616405248642, 0.5808179378509521, 0.3552058041095734),
    @("model_12_1_18", 0.7703788667748346, 0.1992845426132933, 0.3744223106404462, 0.4682590576760349, 0.2541229426860809, 0.1546052843332291, 0.5807619094848633, 0.3551495373249054),
    @("model_12_1_16", 0.7703920433380196, 0.2003054347269134, 0.3742995698190802, 0.4683349829732837, 0.254108339548111, 0.1544081568717957, 0.5808758735656738, 0.3550989031791687),
    @("model_12_1_14", 0.7705046043225413, 0.2012919820957638, 0.3750244399388615, 0.4689601599705492, 0.2539837956428528, 0.1542176753282547, 0.5802029371261597, 0.3546812832355499),
    @("model_12_1_13", 0.7705087615956621, 0.2011550414062432, 0.3750961966283811, 0.4689861725392683, 0.2539791762828827, 0.1542441248893738, 0.580136239528656, 0.3546639382839203),
    @("model_12_1_10", 0.7705501179495036, 0.1996480175017763, 0.3758330709219784, 0.4692374160394654, 0.2539333999156952, 0.154535099864006, 0.5794521570205688, 0.3544961214065552),
    @("model_12_1_8", 0.7705708762758471, 0.1992314851959247, 0.376189836452509, 0.4694070797320247, 0.2539104521274567, 0.154615506529808, 0.5791209936141968, 0.3543828129768372),
    @("model_12_1_9", 0.7706542326660656, 0.2001817733246222, 0.3766032551426509, 0.4698229195005603, 0.2538181841373444, 0.1544320285320282, 0.5787371397018433, 0.3541050553321838),
    @("model_12_1_7", 0.7708329397332174, 0.1996731586598876, 0.3783400532391881, 0.4708813182055462, 0.2536204159259796, 0.1545302271842957, 0.5771248340606689, 0.3533981442451477),
    @("model_12_1_6", 0.7709300596186809, 0.2000670564768837, 0.3791911372247737, 0.4714981031600785, 0.2535129189491272, 0.1544541716575623, 0.5763347148895264, 0.352986216545105),
    @("model_12_1_5", 0.771244044394398, 0.1987109518219667, 0.3823048226352734, 0.4733271939039007, 0.2531654536724091, 0.1547160297632217, 0.5734440684318542, 0.3517645597457886),
    @("model_12_1_4", 0.7717249486766402, 0.200602381289175, 0.3860492302709894, 0.4760659303975757, 0.2526332139968872, 0.1543508172035217, 0.5699679255485535, 0.3499353528022766),
    @("model_12_1_3", 0.772201955086943, 0.2025898400660028, 0.3898793128569239, 0.4788754950543522, 0.2521052956581116, 0.1539670825004578, 0.5664122104644775, 0.3480588793754578),
    @("model_12_1_2", 0.7738168295163852, 0.2091860235607297, 0.4027355969219977, 0.4882942917304769, 0.2503181099891663, 0.1526934504508972, 0.5544769167900085, 0.3417680263519287),
    @("model_12_1_1", 0.774427139253872, 0.2023202297983786, 0.409480314428868, 0.4916551140260499, 0.2496426850557327, 0.1540191322565079, 0.548215389251709, 0.3395233452320099),
    @("model_12_1_0", 0.796495117701392, 0.2888678192338947, 0.6043814120550001, 0.6323867517271473, 0.2252199500799179, 0.1373081803321838, 0.3672768473625183, 0.2455287575721741)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
